# Fixed update to excel issue
#
# 1. Rename "Requested quantity" header to "Weekly_PO_Qty" on the
#    "Weekly Quantity" sheet and "Monthly_PO_Qty" on the "Monthly Trend"
#    sheet.
# 2. Add a new "PO Forecast" worksheet with forecast data (ds, PO_Forecast,
#    yhat_lower, yhat_upper).

$wb = $excel.ActiveWorkbook

# --- 1. Rename headers -----------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the "PO Forecast" sheet -----------------------------------
# Add it directly after the last existing sheet ("Monthly Trend") so it
# lands as the new, final tab without needing a separate Move call.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the header styling already used on the other two sheets (bold,
# thin box border, centered/top aligned) by copying the format from one
# of the existing headers instead of re-building it from scratch.
$wsWeekly.Range("B1").Copy() | Out-Null
$wsForecast.Range("A1:D1").PasteSpecial(-4122) | Out-Null

$data = @(
    @(45347.99999999999, 139, 41.75787721871605, 229.0951921929195),
    @(45375.99999999999, 135, 37.9373939098005, 233.4842256177178),
    @(45389.99999999999, 133, 42.2671626287762, 222.675167047575),
    @(45417.99999999999, 128, 37.84568805866978, 224.532604009343),
    @(45424.99999999999, 127, 28.20484387606845, 224.6365183531958),
    @(45438.99999999999, 125, 31.12804921621887, 218.215109193778),
    @(45445.99999999999, 124, 24.92596709360562, 218.3824988361867),
    @(45459.99999999999, 122, 27.74045055528656, 216.338413079137),
    @(45466.99999999999, 121, 34.38247187785373, 214.4319434384239),
    @(45473.99999999999, 120, 15.56045225609064, 219.0692437525979),
    @(45501.99999999999, 116, 13.40484663658768, 211.5191961912268),
    @(45508.99999999999, 115, 20.23562788176203, 211.7262762475366),
    @(45515.99999999999, 114, 21.20037528387802, 209.3601948933515),
    @(45550.99999999999, 108, 7.836242156043185, 197.6102435718698),
    @(45578.99999999999, 104, 10.64959993899677, 199.6465241403371),
    @(45613.99999999999, 99, 6.812222082459811, 186.4108198385327),
    @(45620.99999999999, 98, 8.435937680249829, 196.0394045208313),
    @(45627.99999999999, 97, -1.346202088689691, 192.8380709788172),
    @(45634.99999999999, 96, 3.918548184673561, 180.8531567268193),
    @(45641.99999999999, 95, -2.880287491599779, 186.6746068516349),
    @(45648.99999999999, 94, -3.289560818895257, 182.7686942704138),
    @(45655.99999999999, 93, -2.748004073465728, 191.1119688639035),
    @(45662.99999999999, 92, -4.509018272833162, 190.6056227688772),
    @(45669.99999999999, 90, -9.347263902633939, 184.8740050265195)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row = $row + 1
}

# The "ds" column holds serial dates - copy the same date-formatted style
# used for column A on the other two sheets.
$wsWeekly.Range("A2").Copy() | Out-Null
$wsForecast.Range("A2:A25").PasteSpecial(-4122) | Out-Null
